# The template contains a single Word field whose field code reads
# " m:enduserdocVar.name " (the "enduserdocVar" portion carries an accent6
# run color). The commit switches the parser from reading real Word fields
# to reading a plain-text "{m:...}" pseudo-field, so the field has to be
# turned into ordinary text runs: "{", "m", ":", "endu", "serdoc", "Var",
# ".name}" - the three middle runs keeping their original color formatting.

$d = $word.ActiveDocument

# Locate the field and the paragraph that hosts it (defensively, instead of
# assuming a fixed paragraph index).
$f = $d.Fields.Item(1)
$fieldStart = $f.Code.Start

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($fieldStart -ge $p.Range.Start -and $fieldStart -lt $p.Range.End) {
        $target = $p
        break
    }
}

$range = $target.Range

# Rebuild the paragraph's content as plain text runs (no more fldChar /
# instrText), preserving the coloured run formatting of "endu", "serdoc"
# and "Var".
$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r><w:t>{</w:t></w:r>
  <w:r><w:t>m</w:t></w:r>
  <w:r><w:t>:</w:t></w:r>
  <w:r>
    <w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr>
    <w:t>endu</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr>
    <w:t>serdoc</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr>
    <w:t>Var</w:t>
  </w:r>
  <w:r><w:t xml:space="preserve">.name}</w:t></w:r>
</w:p>
'@

$range.InsertXML($newParagraphXml)
